$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3 (old rows 3-6 shift down to become rows 4-7)
$ws.Rows.Item(3).Insert()

# Match the date number format used by the date column (D) in the other rows
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat

# Populate the new row 3 with the new weekly record
$ws.Cells.Item(3, 1).Value = 8
$ws.Cells.Item(3, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(3, 3).Value = "Coquimbo"
$ws.Cells.Item(3, 4).Value = 44910
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(3, 6).Value = 100112039
$ws.Cells.Item(3, 7).Value = "Ciboulette"
$ws.Cells.Item(3, 8).Value = "Sin especificar"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 1000
$ws.Cells.Item(3, 11).Value = 1800
$ws.Cells.Item(3, 12).Value = 2000
$ws.Cells.Item(3, 13).Value = 1900
$ws.Cells.Item(3, 14).Value = "`$/docena de atados"
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(3, 16).Value = 633
$ws.Cells.Item(3, 17).Value = 3
$ws.Cells.Item(3, 18).Value = "Hortaliza"
